# Generate Report for Handback
# Updates the timestamp strings recorded for the handback report generation.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: "Latest HO Xliff Generate Date" column (G) ---
# (this timestamp is also mirrored in the de-de sheet's Handoff column, below)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 17:07:52"

# --- "zh-cn" sheet: Correspond Handoff/Handback DateTime columns (H, K) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 17:07:47"
$wsZhCn.Range("K2").Value = "2016-08-24 17:08:24"

# --- "de-de" sheet: Correspond Handoff/Handback DateTime columns (H, K) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-24 17:07:52"
$wsDeDe.Range("K2").Value = "2016-08-24 17:08:32"
